# Auto-generated edit script: apply 2023-10-08 violent crime data update
# Updates column J (year 2023) totals across the Citywide Totals, By Neighborhood,
# and per-neighborhood sheets, per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('J2').Value = 5903
$ws.Range('J3').Value = 6317
$ws.Range('J4').Value = 1368
$ws.Range('J5').Value = 483
$ws.Range('J6').Value = 8108
$ws.Range('J7').Value = 22179

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('J6').Value = 203
$ws.Range('J7').Value = 308

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('J2').Value = 381
$ws.Range('J3').Value = 425
$ws.Range('J6').Value = 472
$ws.Range('J7').Value = 1391

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('J3').Value = 166
$ws.Range('J6').Value = 122
$ws.Range('J7').Value = 450

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('J4').Value = 43
$ws.Range('J6').Value = 356
$ws.Range('J7').Value = 1030

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('J3').Value = 114
$ws.Range('J7').Value = 324

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('J3').Value = 232
$ws.Range('J7').Value = 680

$ws = $wb.Worksheets.Item('New City')
$ws.Range('J3').Value = 160
$ws.Range('J7').Value = 557

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('J6').Value = 90
$ws.Range('J7').Value = 348

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('J2').Value = 176
$ws.Range('J6').Value = 162
$ws.Range('J7').Value = 658
$ws.Range('J8').Value = 1391
$ws.Range('J10').Value = 155
$ws.Range('J11').Value = 350
$ws.Range('J12').Value = 44
$ws.Range('J13').Value = 28
$ws.Range('J14').Value = 110
$ws.Range('J16').Value = 88
$ws.Range('J19').Value = 659
$ws.Range('J20').Value = 460
$ws.Range('J22').Value = 56
$ws.Range('J24').Value = 70
$ws.Range('J29').Value = 1224
$ws.Range('J33').Value = 1030
$ws.Range('J36').Value = 307
$ws.Range('J37').Value = 680
$ws.Range('J39').Value = 12
$ws.Range('J42').Value = 933
$ws.Range('J48').Value = 260
$ws.Range('J51').Value = 273
$ws.Range('J52').Value = 559
$ws.Range('J53').Value = 308
$ws.Range('J54').Value = 434
$ws.Range('J57').Value = 99
$ws.Range('J60').Value = 131
$ws.Range('J63').Value = 80
$ws.Range('J65').Value = 557
$ws.Range('J67').Value = 840
$ws.Range('J71').Value = 75
$ws.Range('J73').Value = 213
$ws.Range('J76').Value = 339
$ws.Range('J77').Value = 166
$ws.Range('J79').Value = 634
$ws.Range('J83').Value = 450
$ws.Range('J84').Value = 188
$ws.Range('J85').Value = 911
$ws.Range('J86').Value = 142
$ws.Range('J89').Value = 290
$ws.Range('J91').Value = 252
$ws.Range('J92').Value = 69
$ws.Range('J95').Value = 324
$ws.Range('J99').Value = 348
$ws.Range('J101').Value = 22179

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('J2').Value = 210
$ws.Range('J3').Value = 317
$ws.Range('J6').Value = 225
$ws.Range('J7').Value = 840

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('J6').Value = 57
$ws.Range('J7').Value = 188

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('J2').Value = 102
$ws.Range('J6').Value = 210
$ws.Range('J7').Value = 434

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('J2').Value = 373
$ws.Range('J3').Value = 429
$ws.Range('J7').Value = 1224

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('J6').Value = 130
$ws.Range('J7').Value = 260

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('J2').Value = 158
$ws.Range('J3').Value = 192
$ws.Range('J5').Value = 24
$ws.Range('J7').Value = 659

$ws = $wb.Worksheets.Item('River North')
$ws.Range('J3').Value = 68
$ws.Range('J6').Value = 190
$ws.Range('J7').Value = 339

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('J3').Value = 22
$ws.Range('J6').Value = 39
$ws.Range('J7').Value = 110

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('J2').Value = 46
$ws.Range('J7').Value = 162

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('J3').Value = 191
$ws.Range('J6').Value = 482
$ws.Range('J7').Value = 933

$ws = $wb.Worksheets.Item('Boystown')
$ws.Range('J4').Value = 5
$ws.Range('J6').Value = 28

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('J6').Value = 82
$ws.Range('J7').Value = 155

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range('J6').Value = 17
$ws.Range('J7').Value = 70

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('J2').Value = 70
$ws.Range('J7').Value = 252

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('J2').Value = 178
$ws.Range('J3').Value = 217
$ws.Range('J7').Value = 634

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('J2').Value = 126
$ws.Range('J3').Value = 160
$ws.Range('J4').Value = 41
$ws.Range('J7').Value = 460

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('J3').Value = 100
$ws.Range('J7').Value = 307

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('J2').Value = 201
$ws.Range('J3').Value = 198
$ws.Range('J7').Value = 658

$ws = $wb.Worksheets.Item('Greektown')
$ws.Range('J4').Value = 2
$ws.Range('J6').Value = 12

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('J6').Value = 147
$ws.Range('J7').Value = 350

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('J3').Value = 55
$ws.Range('J7').Value = 213

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('J6').Value = 67
$ws.Range('J7').Value = 176

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range('J2').Value = 19
$ws.Range('J7').Value = 69

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('J3').Value = 81
$ws.Range('J7').Value = 290

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('J3').Value = 31
$ws.Range('J6').Value = 45

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('J4').Value = 76
$ws.Range('J7').Value = 142

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('J2').Value = 65
$ws.Range('J7').Value = 273

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range('J6').Value = 42
$ws.Range('J7').Value = 99

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('J2').Value = 47
$ws.Range('J7').Value = 131

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('J5').Value = 19
$ws.Range('J7').Value = 911

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range('J2').Value = 26
$ws.Range('J7').Value = 56

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range('J3').Value = 22
$ws.Range('J7').Value = 75

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('J3').Value = 56
$ws.Range('J7').Value = 166

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('J2').Value = 134
$ws.Range('J6').Value = 227
$ws.Range('J7').Value = 559

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range('J2').Value = 6
$ws.Range('J7').Value = 44

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range('J3').Value = 7
$ws.Range('J7').Value = 88

Write-Host "Updated 157 cells across 46 sheets"